$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Add the new "Sheet1" worksheet after the existing "loginData" sheet
# ---------------------------------------------------------------------
$loginData = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $loginData)
$newSheet.Name = "Sheet1"

# ---------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------
$headers = @("Batch name","class topic","class description","Select class dates","No of classes","staff name","status")
for ($i = 0; $i -lt $headers.Length; $i++) {
  $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------
# Data rows (A..H). Column D holds a date serial formatted as a date;
# column E holds a plain number; everything else is text.
# ---------------------------------------------------------------------
$rows = @(
  @{ A="Micro service -03"; B="core Java";   C="java";         D=45711; E=3; F="saranya m";        G="Active";   H="with valid data" },
  @{ A="Micro service-01";  B="IK0en";       C='""';           D=45718; E=2; F="123InvalidTopic!";  G="InActive"; H="with invalid data" },
  @{ A="SMPO-0001";         B="new test 7";  C=$null;          D=45712; E=1; F="Kevin Thomas";      G="Active";   H="with mandatory fields" },
  @{ A="Python101";         B="jmeter123";   C="HTML";         D=46441; E=1; F="Geetha takur";      G="Active";   H="with ordinary fields" },
  @{ A="Micro service-01";  B="Vidhya Test"; C="Vidhya desc";  D=45712; E=1; F="Getha  Takur";      G="Active";   H="close cancel button" }
)

$r = 2
foreach ($row in $rows) {
  $newSheet.Cells.Item($r, 1).Value = $row.A
  $newSheet.Cells.Item($r, 2).Value = $row.B
  if ($row.C -ne $null) {
    $newSheet.Cells.Item($r, 3).Value = $row.C
  }
  $newSheet.Cells.Item($r, 4).Value = $row.D
  $newSheet.Cells.Item($r, 5).Value = $row.E
  $newSheet.Cells.Item($r, 6).Value = $row.F
  $newSheet.Cells.Item($r, 7).Value = $row.G
  $newSheet.Cells.Item($r, 8).Value = $row.H
  $r++
}

# Apply the date number format once, then copy/paste-special (formats
# only) onto the remaining date cells so they all share one style index
# instead of each getting its own duplicate cellXfs entry.
$newSheet.Range("D2").NumberFormat = "mm-dd-yy"
$newSheet.Range("D2").Copy()
$newSheet.Range("D3:D6").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Column widths (best-effort match of the authored widths)
# ---------------------------------------------------------------------
$newSheet.Columns.Item(2).ColumnWidth = 19.92
$newSheet.Columns.Item(3).ColumnWidth = 16.92
$newSheet.Columns.Item(4).ColumnWidth = 18.09
$newSheet.Columns.Item(6).ColumnWidth = 18.26
$newSheet.Columns.Item(8).ColumnWidth = 20.26

# ---------------------------------------------------------------------
# Selection state: loginData no longer the active/selected sheet;
# new Sheet1 becomes the active tab with D9/E9 selections respectively.
# ---------------------------------------------------------------------
$loginData.Range("D9").Select() | Out-Null
$newSheet.Range("E9").Select() | Out-Null
